$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" data column (Q) is being added after the existing 2019 column (P).
# Column O already carries the exact per-row cell-style pattern that column P
# uses (plain style on D:K/M, "applyFill" style on L/N/O/P), so copying column O
# into a freshly inserted column Q reproduces the same style indices that the
# real edit ended up with for column Q, instead of Excel/iron_native minting a
# brand new style record.
$null = $ws.Range("O1:O14").Copy()
$null = $ws.Range("Q1:Q14").Insert(-4161)

# Header + the 2020 values for each indicator row
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 0.02
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("Q9").Value = 0.54
$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("Q14").Value = 0

# Match the selection recorded in the sheet view
$null = $ws.Range("N19").Select()
